$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTickers = @("TAO-USD", "IMX-USD", "GRT-USD", "PEPE-USD", "MNT-USD")
$startRow = 424

for ($i = 0; $i -lt $newTickers.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newTickers[$i]
}
